# Adds a second student record (row 3) to the import-siswa sheet, mirroring
# the layout of row 2: No, Username, Email (hyperlinked), Password.
#
# Shared strings must come out in this order: "rdfd", "TrianNurizkillah6",
# "triannurizkillah6@gmail.com" -> so we write column D before B before C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 2
$ws.Range("D3").Value = "rdfd"
$ws.Range("B3").Value = "TrianNurizkillah6"
$ws.Range("C3").Value = "triannurizkillah6@gmail.com"

# Turn the new email cell into a mailto: hyperlink, then make sure it carries
# the same "Hyperlink" cell style as the existing C2 link.
[void]$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:triannurizkillah6@gmail.com")
$ws.Range("C3").Style = "Hyperlink"

# Match the workbook's on-disk selection state (active cell moves to C3).
[void]$ws.Range("C3").Select()
